# Update the supervisor upload list: replace placeholder supervisor names
# with the real supervisor list and append four new supervisor rows
# (tid00011-tid00014), per "modify on recieving student credit inputs".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-11: supervisor name (col A) and priority (col D) ---
# tid/password columns (B/C) are unchanged for these rows.

$ws.Range("A2").Value = "Dr LIU Yang"
$ws.Range("D2").Value = 2

$ws.Range("A3").Value = "Dr. ZHANG, Eric Lu"
$ws.Range("D3").Value = 2

$ws.Range("A4").Value = "Dr. WAN, Renjie"
$ws.Range("D4").Value = 3

$ws.Range("A5").Value = "Dr Yu, Wilson Shih Bun"
$ws.Range("D5").Value = 3

$ws.Range("A6").Value = "Prof. CHEUNG, Yiu Ming"
$ws.Range("D6").Value = 2

$ws.Range("A7").Value = "Dr. ZHOU, Kaiyang"
$ws.Range("D7").Value = 3

$ws.Range("A8").Value = "Prof. Xu, Jianliang"
$ws.Range("D8").Value = 1

$ws.Range("A9").Value = "Dr. HAN, Bo"
$ws.Range("D9").Value = 2

$ws.Range("A10").Value = "Dr. DAI, Henry Hong Ning"
$ws.Range("D10").Value = 3

$ws.Range("A11").Value = "Prof. YUEN, Pong Chi"
$ws.Range("D11").Value = 2

# --- Append new rows 12-15 for newly onboarded supervisors ---

$ws.Range("A12").Value = "Prof. LEUNG,Yiu Wing"
$ws.Range("B12").Value = "tid00011"
$ws.Range("C12").Value = "tpw00011"
$ws.Range("D12").Value = 2

$ws.Range("A13").Value = "Dr. FENG, Jian"
$ws.Range("B13").Value = "tid00012"
$ws.Range("C13").Value = "tpw00012"
$ws.Range("D13").Value = 3

$ws.Range("A14").Value = "Dr. WANG, Juncheng"
$ws.Range("B14").Value = "tid00013"
$ws.Range("C14").Value = "tpw00013"
$ws.Range("D14").Value = 3

$ws.Range("A15").Value = "Dr. YANG, Renchi"
$ws.Range("B15").Value = "tid00014"
$ws.Range("C15").Value = "tpw00014"
$ws.Range("D15").Value = 3

# --- Update the view: scrolled down to the new rows, selection on D16 ---
$null = $ws.Range("D16").Select()

$wb.Save()
